# Auto-generated script applying the scheduled-runner profit updates
# to the Sargatanas_Profits workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 12500308
$ws.Range("I41").Value = 12500308
$ws.Range("K41").Value = 12500308
$ws.Range("M41").Value = -12499868
$ws.Range("H62").Value = 71450170
$ws.Range("I62").Value = 142860660
$ws.Range("J62").Value = 39685.57
$ws.Range("K62").Value = 142860660
$ws.Range("L62").Value = 39685.57
$ws.Range("M62").Value = -142860036
$ws.Range("N62").Value = -40933.57
$ws.Range("H65").Value = 71450170
$ws.Range("I65").Value = 142860660
$ws.Range("J65").Value = 39685.57
$ws.Range("K65").Value = 714303300
$ws.Range("L65").Value = 198427.85
$ws.Range("M65").Value = -714300180
$ws.Range("N65").Value = -204667.85
$ws.Range("H106").Value = 250002340
$ws.Range("I106").Value = 250002340
$ws.Range("K106").Value = 250002340
$ws.Range("M106").Value = -250001709
$ws.Range("H107").Value = 93752730
$ws.Range("I107").Value = 62501950
$ws.Range("J107").Value = 125003500
$ws.Range("K107").Value = 62501950
$ws.Range("L107").Value = 125003500
$ws.Range("M107").Value = -62500030
$ws.Range("N107").Value = -125007340
$ws.Range("H113").Value = 125015870
$ws.Range("J113").Value = 125015870
$ws.Range("L113").Value = 125015870
$ws.Range("N113").Value = -125022378
$ws.Range("H116").Value = 10874503
$ws.Range("I116").Value = 31252562
$ws.Range("K116").Value = 31252562
$ws.Range("M116").Value = -31249120
$ws.Range("H118").Value = 751
$ws.Range("I118").Value = 788.2857
$ws.Range("K118").Value = 2364.8571
$ws.Range("M118").Value = -707.8571000000002
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 1984.8939
$ws.Range("I132").Value = 1933.4603
$ws.Range("K132").Value = 5800.3809
$ws.Range("M132").Value = -3270.3809
$ws.Range("H134").Value = 80279.5
$ws.Range("J134").Value = 80279.5
$ws.Range("L134").Value = 80279.5
$ws.Range("N134").Value = -90419.5
$ws.Range("H138").Value = 5902.82
$ws.Range("I138").Value = 2599.6428
$ws.Range("J138").Value = 7187.3887
$ws.Range("K138").Value = 7798.928400000001
$ws.Range("L138").Value = 21562.1661
$ws.Range("M138").Value = -2658.928400000001
$ws.Range("N138").Value = -31842.1661

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 10376
$ws.Range("J43").Value = 10376
$ws.Range("L43").Value = 10376
$ws.Range("N43").Value = -11002
$ws.Range("H110").Value = 333333340
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 333333340
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 333333340
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -333337430
$ws.Range("H132").Value = 4624.098
$ws.Range("I132").Value = 2092.0645
$ws.Range("J132").Value = 8548.75
$ws.Range("K132").Value = 6276.193499999999
$ws.Range("L132").Value = 25646.25
$ws.Range("M132").Value = -3746.193499999999
$ws.Range("N132").Value = -30706.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3068.4092
$ws.Range("I105").Value = 2230.7693
$ws.Range("K105").Value = 2230.7693
$ws.Range("M105").Value = -483.7692999999999
$ws.Range("H134").Value = 6643.7715
$ws.Range("I134").Value = 3192.2
$ws.Range("K134").Value = 9576.599999999999
$ws.Range("M134").Value = -7041.599999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 41336144
$ws.Range("I86").Value = 20223906
$ws.Range("J86").Value = 66670830
$ws.Range("K86").Value = 20223906
$ws.Range("L86").Value = 66670830
$ws.Range("M86").Value = -20222783
$ws.Range("N86").Value = -66673076
$ws.Range("H89").Value = 41336144
$ws.Range("I89").Value = 20223906
$ws.Range("J89").Value = 66670830
$ws.Range("K89").Value = 101119530
$ws.Range("L89").Value = 333354150
$ws.Range("M89").Value = -101113914
$ws.Range("N89").Value = -333365382
$ws.Range("H99").Value = 10001.5625
$ws.Range("J99").Value = 7793.1816
$ws.Range("L99").Value = 7793.1816
$ws.Range("N99").Value = -10789.1816
$ws.Range("H126").Value = 10001.5625
$ws.Range("J126").Value = 7793.1816
$ws.Range("L126").Value = 23379.5448
$ws.Range("N126").Value = -28319.5448
$ws.Range("H132").Value = 6322.45
$ws.Range("I132").Value = 2214
$ws.Range("K132").Value = 6642
$ws.Range("M132").Value = -4112

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 125894.19
$ws.Range("I2").Value = 373.63635
$ws.Range("K2").Value = 2241.8181
$ws.Range("M2").Value = -2128.8181

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1317.9269
$ws.Range("I97").Value = 1200.0385
$ws.Range("J97").Value = 1522.2667
$ws.Range("K97").Value = 1200.0385
$ws.Range("L97").Value = 1522.2667
$ws.Range("M97").Value = -704.0385000000001
$ws.Range("N97").Value = -2514.2667
$ws.Range("H107").Value = 1334408.1
$ws.Range("I107").Value = 1600989.8
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 1600989.8
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = -1599069.8
$ws.Range("N107").Value = -5340
$ws.Range("H122").Value = 52980.6
$ws.Range("I122").Value = 73229.78999999999
$ws.Range("K122").Value = 219689.37
$ws.Range("M122").Value = -217239.37
$ws.Range("H132").Value = 7902.357
$ws.Range("I132").Value = 2910
$ws.Range("J132").Value = 16888.6
$ws.Range("K132").Value = 8730
$ws.Range("L132").Value = 50665.8
$ws.Range("M132").Value = -6200
$ws.Range("N132").Value = -55725.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4682.8945
$ws.Range("I122").Value = 3281.1667
$ws.Range("J122").Value = 5329.846
$ws.Range("K122").Value = 9843.500100000001
$ws.Range("L122").Value = 15989.538
$ws.Range("M122").Value = -7393.500100000001
$ws.Range("N122").Value = -20889.538
$ws.Range("H132").Value = 16674982
$ws.Range("I132").Value = 35719284
$ws.Range("K132").Value = 107157852
$ws.Range("M132").Value = -107155322
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 25204432
$ws.Range("I122").Value = 45821530
$ws.Range("J122").Value = 5760.3335
$ws.Range("K122").Value = 137464590
$ws.Range("L122").Value = 17281.0005
$ws.Range("M122").Value = -137462140
$ws.Range("N122").Value = -22181.0005
$ws.Range("H126").Value = 142863420
$ws.Range("I126").Value = 250005250
$ws.Range("K126").Value = 750015750
$ws.Range("M126").Value = -750013280
$ws.Range("H132").Value = 21768376
$ws.Range("I132").Value = 31257946
$ws.Range("J132").Value = 77928.57000000001
$ws.Range("K132").Value = 93773838
$ws.Range("L132").Value = 233785.71
$ws.Range("M132").Value = -93771308
$ws.Range("N132").Value = -238845.71
$ws.Range("H135").Value = 85238.336
$ws.Range("J135").Value = 85238.336
$ws.Range("L135").Value = 85238.336

Write-Host "Applied $(($wb.Worksheets | Measure-Object).Count) worksheet(s); profit figures refreshed."
